$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 20: refined mean_Intake / sem_Intake values
$ws.Range("P20").Value = 59.82145489760461
$ws.Range("Q20").Value = 56.969700048919947

# Row 54: previously placeholder 65535 / blank -> real computed values
$ws.Range("P54").Value = 56.220744202549618
$ws.Range("Q54").Value = 19.425478785356162

# Row 56: previously placeholder 65535 / blank -> real computed values
$ws.Range("P56").Value = 114.90023813898675
$ws.Range("Q56").Value = 50.044468626942816

# Row 59: previously placeholder 65535 / blank -> real computed values
$ws.Range("P59").Value = 61.087244108135195
$ws.Range("Q59").Value = 23.644092541978718

# Row 64: refined mean_Intake / sem_Intake values
$ws.Range("P64").Value = 42.385652341041457
$ws.Range("Q64").Value = 9.7041622791667983

# Row 65: previously placeholder 65535 / blank -> real computed values
$ws.Range("P65").Value = 37.897737798460241
$ws.Range("Q65").Value = 17.144928103836946
